# Jit Model Class Design.pptx -- "Stage Engine Linkage (Step1)"
#
# 1) The footer "datetimeFigureOut" date placeholder (present on the slide
#    master and every slide layout) is refreshed from 2020/1/30 to 2020/2/5.
# 2) On slide 1, the rounded-rectangle shape that just reads "Stage" is
#    relabeled "StageSubset：X" (matching the sibling "Stage : Y" / "Stage : Z"
#    boxes), split into three runs so the full-width colon keeps its own
#    run, same as the other boxes.

$p = $ppt.ActivePresentation

$oldDate = "2020/1/30"
$newDate = "2020/2/5"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }
        $hasPh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) { $hasPh = $true }
        } catch {
            $hasPh = $false
        }
        if (-not $hasPh) { continue }
        if (-not $shp.TextFrame.HasText) { continue }
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master

# Every custom (slide) layout hanging off the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Update-DatePlaceholder $master.CustomLayouts.Item($li)
}

# --- Slide 1: "Stage" -> "StageSubset：X" -------------------------------

$slide = $p.Slides.Item(1)
$stageShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "Stage") {
                $stageShape = $shp
            }
        }
    }
}

if ($stageShape -ne $null) {
    $tr = $stageShape.TextFrame.TextRange
    $tr.Text = "StageSubsetXX"

    # First run keeps "StageSubset" as-is (already carries the original
    # bold/underline/accent4 formatting). Re-touch the trailing two
    # characters so the engine splits them into independent runs, then
    # fix up their actual text.
    $colonRun = $tr.Characters(12, 1)
    $colonRun.Font.Bold = $true
    $colonRun.Text = "："

    $xRun = $tr.Characters(13, 1)
    $xRun.Font.Bold = $true
    $xRun.Text = "X"
}
